$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '61.547.57'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.447.56'
$ws.Range("E3").Value = '  +2.17%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '579.28'
$ws.Range("E5").Value = '  +1.27%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '144.49'
$ws.Range("E6").Value = '  +6.02%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.450.15'
$ws.Range("E7").Value = '  +2.29%  '
$ws.Range("E8").Value = '  -0.03%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.476'
$ws.Range("E9").Value = '  +1.75%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.62'
$ws.Range("E10").Value = '  +0.19%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.125'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.388'
$ws.Range("E12").Value = '  +2.48%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.036.01'
$ws.Range("E13").Value = '  +2.14%  '
$ws.Range("E14").Value = '  +9.17%  '
$ws.Range("E15").Value = '  -0.84%  '
$ws.Range("B16").Value = 'ShibaInu'
$ws.Range("C16").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000173'
$ws.Range("E16").Value = '  +1.79%  '
$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.453.73'
$ws.Range("E17").Value = '  +2.38%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '61.699.31'
$ws.Range("E18").Value = '  +0.95%  '
$ws.Range("E19").Value = '  +9.06%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.27'
$ws.Range("E20").Value = '  +3.67%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.54'
$ws.Range("E21").Value = '  +2.63%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '389.03'
$ws.Range("E22").Value = '  +3.48%  '
$ws.Range("E23").Value = '  +3.44%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '73.46'
$ws.Range("E24").Value = '  +3.45%  '
$ws.Range("E25").Value = '  +0.15%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("E26").Value = '  +0.02%  '
$ws.Range("E27").Value = '  +0.34%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '3.590.98'
$ws.Range("E28").Value = '  +2.09%  '
$ws.Range("E29").Value = '  +2.70%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.60'
$ws.Range("E30").Value = '  +3.18%  '
$ws.Range("E31").Value = '  +0.16%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.16'
$ws.Range("E32").Value = '  +1.65%  '
$ws.Range("E33").Value = '  +2.45%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.45'
$ws.Range("E34").Value = '  -11.27%  '
$ws.Range("E35").Value = '  -0.02%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '24.03'
$ws.Range("E36").Value = '  +2.95%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.475.55'
$ws.Range("E37").Value = '  +2.34%  '
$ws.Range("E38").Value = '  +3.28%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.12'
$ws.Range("E39").Value = '  +0.38%  '
$ws.Range("E40").Value = '  +0.72%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '166.70'
$ws.Range("E41").Value = '  +1.24%  '
$ws.Range("E42").Value = '  +3.26%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '27.86'
$ws.Range("E43").Value = '  +11.03%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.804'
$ws.Range("E44").Value = '  +3.84%  '
$ws.Range("E45").Value = '  +0.05%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '42.40'
$ws.Range("E46").Value = '  +1.87%  '
$ws.Range("E47").Value = '  +4.09%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.73'
$ws.Range("E48").Value = '  +3.00%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.569.07'
$ws.Range("E49").Value = '  +1.33%  '
$ws.Range("E50").Value = '  -1.53%  '
$ws.Range("E51").Value = '  +2.60%  '
